$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(32, 45693.18645833333, 8, 6, 195, 380, 356, 384, 2664, 384, 1216, 119, 304, 30, 2985, 3873),
    @(33, 45693.191284722219, 8, 6, 195, 380, 356, 384, 2664, 384, 1216, 119, 304, 30, 2985, 3873),
    @(34, 45694.214224537034, 8, 6, 196, 381, 357, 384, 2664, 384, 1216, 119, 304, 30, 2993, 3903),
    @(35, 45694.21837962963, 8, 6, 196, 381, 357, 384, 2664, 384, 1216, 119, 304, 30, 2993, 3903),
    @(36, 45694.224363425928, 8, 6, 196, 381, 357, 384, 2664, 384, 1216, 119, 304, 30, 2993, 3903),
    @(37, 45694.928773148145, 8, 6, 198, 381, 357, 384, 2664, 384, 1216, 119, 304, 30, 2999, 3934)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($col = 1; $col -le 15; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
